$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# 102-BULLS-SHORTER (row 4) "Not Same Time As" constraint: shift the no-go window an hour later
$ws.Range("H4").Value = "no 5:40pm-6:40pm"

# 106-PATRIOTS (row 8) "Not Same Time As" constraint: used to duplicate BULLS-SHORTER's old
# window, give it its own distinct (later) time-based constraint
$ws.Range("H8").Value = "no 6:40pm-7:40pm"

# 107-SHOOTIN TIGERS-JENNINGS (row 9) "Not Same Time As" constraint: change category from
# "only" to "pref" for the day-based constraint
$ws.Range("H9").Value = "pref thursday, no saturday"

# Update the active selection on the sheet
$ws.Range("H18").Select()
